$wb = $excel.ActiveWorkbook

# --- Sheet "run_settings" ---
$ws1 = $wb.Worksheets.Item("run_settings")

# B3 (uniprot_list): value changed from 2 to 94
$ws1.Range("B3").Value = 94

# Several run_* / slice_* flags flipped from TRUE to FALSE
$ws1.Range("B19").Value = "FALSE"   # run_create_csv_from_uniprot_flatfile
$ws1.Range("B20").Value = "FALSE"   # run_setup_df_file_locations
$ws1.Range("B24").Value = "FALSE"   # run_parse_simap_to_csv
$ws1.Range("B27").Value = "FALSE"   # slice_TMDs_from_homologues
$ws1.Range("B29").Value = "FALSE"   # run_calculate_AAIMON_ratios

# Update the saved selection for this sheet
$ws1.Range("B30").Select()

# --- Sheet "file_locations" ---
$ws2 = $wb.Worksheets.Item("file_locations")
$ws2.Activate()
$ws2.Range("A21").Select()

# --- Sheet "variables" ---
$ws3 = $wb.Worksheets.Item("variables")
$ws3.Activate()
$ws3.Range("B33").Select()

# run_settings stays the sheet that is active/selected when the file is saved
$ws1.Activate()
